$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add "Username" / "Password" columns, matching the
# bold header style already used by A1/B1 ---
$ws.Range("C1").Value = "Username"
$ws.Range("D1").Value = "Password"
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 2: new Username/Password sample values, plain (default) style ---
$ws.Range("C2").Value = "standard_user"
$ws.Range("D2").Value = "secret_sauce"

# --- Row 3: a second "Valid URL" hyperlink row plus sample creds ---
$ws.Range("A3").Value = "https://www.saucedemo.com/inventory.html"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.saucedemo.com/inventory.html", "", "", "https://www.saucedemo.com/inventory.html")
$ws.Range("C3").Value = "abcdefghpoiuyt"
$ws.Range("D3").Value = "qwertyuiop"

# Re-apply the existing hyperlink look (blue font) to A3, matching A2/B2 --
# Hyperlinks.Add resets formatting to a generic style, so copy it after.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)      # xlPasteFormats

# Give row 3 the same custom-height treatment as rows 1/2
$ws.Rows.Item(3).RowHeight = 89.55

# --- Column widths for the two new columns, and a tiny bump to column B ---
$ws.Columns.Item(2).ColumnWidth = 36
$ws.Columns.Item(3).ColumnWidth = 19.6
$ws.Columns.Item(4).ColumnWidth = 22.6

# --- Selection follows the last entered cell ---
$ws.Range("D3").Select()
